$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$korlash = "('Korlash, Heir to Blackblade', ['{2}{B}{B}', 'Legendary Creature " + [char]0x2014 + " Zombie Warrior', 'Korlash, Heir to Blackblade" + [char]0x2019 + "s power and toughness are each equal to the number of Swamps you control.', '{1}{B}: Regenerate Korlash.', 'Grandeur " + [char]0x2014 + " Discard another card named Korlash, Heir to Blackblade: Search your library for up to two Swamp cards, put them onto the battlefield tapped, then shuffle your library.', '*/*'])"

$storm = "('Storm Entity', ['{1}{R}', 'Creature " + [char]0x2014 + " Elemental', 'Haste', 'Storm Entity enters the battlefield with a +1/+1 counter on it for each other spell cast this turn.', '1/1'])"

# Clear rows 4-14 first (contents only, so the range doesn't shift and
# the used range / dimension shrinks to A1:A3 once A2/A3 are rewritten).
$ws.Range("A4:A14").ClearContents() | Out-Null

$ws.Range("A2").Value = $korlash
$ws.Range("A3").Value = $storm
